# Apply hybrid bold + color highlighting to quantitative metrics in
# specific bullet/paragraph lines of the resume.
#
# Strategy: for each target paragraph (identified by its exact, full
# pre-edit text), re-find it by exact match, then within that
# paragraph's Range, locate each metric substring IN ORDER (searching
# strictly after the end of the previously-found metric so repeated
# tokens - e.g. two "87%" style matches in different paragraphs, or a
# token that could appear twice - are handled correctly) and set
# Font.Bold + Font.Color on just that sub-range. Word's editor
# automatically splits/creates the surrounding runs, which is exactly
# the <w:r> split seen in the target diff.

$d = $word.ActiveDocument

# Word (VBA/COM) colors are BGR-packed integers (RGB() macro order),
# not the RRGGBB hex used in OOXML w:color. Convert once.
function HexToWordColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$metricColor = HexToWordColor "2C3E50"

# Highlights a list of literal substrings, in left-to-right order, inside
# the given paragraph's range, applying bold + the metric color to each.
function Highlight-Metrics($paragraph, [string[]]$metrics) {
    $cursor = $paragraph.Range.Start
    $paraEnd = $paragraph.Range.End
    foreach ($metric in $metrics) {
        $searchRng = $d.Range($cursor, $paraEnd)
        $found = $searchRng.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $searchRng.Font.Bold = $true
            $searchRng.Font.Color = $metricColor
            $cursor = $searchRng.End
        }
    }
}

# Finds the single paragraph whose trimmed text matches $exactText
# exactly, and runs $action against it. $exactText should be given
# without the leading bullet glyph / bullet+space, matched via
# Contains/EndsWith so we don't need to fuss over the bullet char.
function Find-ParagraphAndHighlight([string]$containsText, [string]$endsWithText, [string[]]$metrics) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Contains($containsText) -and $t.TrimEnd().EndsWith($endsWithText)) {
            Highlight-Metrics $p $metrics
            return
        }
    }
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
Find-ParagraphAndHighlight `
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" `
    "from 23% to 64%" `
    @("23%", "64%")

# 2) "Achieved 87% prediction accuracy ... margins from ±4.2% to ±2.1%" (long form)
Find-ParagraphAndHighlight `
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from" `
    "±2.1%" `
    @("87%", "71%", "±4.2%", "±2.1%")

# 3) "Wrote RFP and analyzed bids from 1,200 vendors ..."
Find-ParagraphAndHighlight `
    "Wrote RFP and analyzed bids from 1,200 vendors for research platform development" `
    "research platform development" `
    @("1,200")

# 4) "Created comprehensive meta-analysis framework ... $400M ... $1B+"
Find-ParagraphAndHighlight `
    "Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" `
    "now valued at `$1B+" `
    @("`$400M", "`$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Find-ParagraphAndHighlight `
    "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" `
    "organizations `$4.7M" `
    @("73.5%", "`$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form)
Find-ParagraphAndHighlight `
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" `
    "industry standard of 71%" `
    @("87%", "71%")
